# "Design Approach section finished"
#
# Replaces the placeholder sentence in the "Design Approach" paragraph with
# the finished write-up, and relocates Word's "_GoBack" (last-edit-location)
# bookmark from its previous spot (end of the Testing/Evaluation intro
# paragraph) into the middle of the newly-written text, matching where the
# author's cursor was when they stopped typing.

$d = $word.ActiveDocument

$newText = "As this is a large project with five stages and several components per stage, we took a multi-step approach to design. First, we needed a complete block diagram detailing all the components needed for each stage, as in Figure 1. After verifying with the course material and online sources that we had all the necessary components and control signals, we wrote the descriptions of these components stage by stage. As we were a group of four, we could do the component design in parallel. Upon completion of the components of a stage, they were integrated together as a single block in a higher-level module. This module contained all the input and output signals of its respective stage and connected the ports of each component to the necessary signals. We should note that the latches dividing each of the five stages were written as modules of their own. Finally, after each stage had successfully integrated the components within them, a top-level module was designed for connecting the entire pipeline together. In this top-level module, inter-stage signals such as forwarding signals, write enabling and write back signals were implemented and connected to their required ports."

# 1. Replace the old placeholder sentence with the finished paragraph text.
$d.Content.Find.Execute("The pipelined processor was designed to…", $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# 2. Work out where "the component " ends / "design in parallel." begins
#    inside the freshly-inserted text, so we can drop the bookmark there.
$locate = $d.Content
$locate.Find.Execute("the component design in parallel.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkPos = $d.Range($locate.Start + 14, $locate.Start + 14)

# 3. Re-adding a bookmark named "_GoBack" moves it (Word keeps only one,
#    silently dropping the previous instance wherever it was).
$d.Bookmarks.Add("_GoBack", $bookmarkPos)
